# Applies crypto price/volume updates per commit "Updated cryptos list on Sun Jun 30 17:42:44 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.762.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.412.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.995.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.409.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.812.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("E20").Value = "  +2.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.185"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.447.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.90%  "
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.786"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.500.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  -3.56%  "
$ws.Range("E51").Value = "  -0.08%  "
